$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.723.72"
$ws.Range("E2").Value = "  +2.17%  "

$ws.Range("D3").Value = "'3.544.61"
$ws.Range("E3").Value = "  +1.05%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").Value = "'608.66"
$ws.Range("E5").Value = "  +4.51%  "

$ws.Range("D6").Value = "'174.38"
$ws.Range("E6").Value = "  +0.60%  "

$ws.Range("D7").Value = "'0.619"
$ws.Range("E7").Value = "  -0.34%  "

$ws.Range("D8").Value = "'3.536.47"
$ws.Range("E8").Value = "  +1.01%  "

$ws.Range("E9").Value = "  -0.12%  "

$ws.Range("D10").Value = "'0.201"
$ws.Range("E10").Value = "  +6.31%  "

$ws.Range("D11").Value = "'6.76"
$ws.Range("E11").Value = "  +0.26%  "

$ws.Range("D12").Value = "'0.587"
$ws.Range("E12").Value = "  -1.39%  "

$ws.Range("D13").Value = "'47.63"
$ws.Range("E13").Value = "  +1.62%  "

$ws.Range("D14").Value = "'0.0000281"
$ws.Range("E14").Value = "  +1.71%  "

$ws.Range("D15").Value = "'4.106.62"
$ws.Range("E15").Value = "  +0.69%  "

$ws.Range("D16").Value = "'627.37"
$ws.Range("E16").Value = "  -7.33%  "

$ws.Range("D17").Value = "'8.45"
$ws.Range("E17").Value = "  -2.98%  "

$ws.Range("D18").Value = "'70.653.91"
$ws.Range("E18").Value = "  +1.97%  "

$ws.Range("D19").Value = "'3.541.26"
$ws.Range("E19").Value = "  +0.65%  "

$ws.Range("E20").Value = "  -1.73%  "

$ws.Range("D21").Value = "'17.47"
$ws.Range("E21").Value = "  +0.09%  "

$ws.Range("D22").Value = "'10.03"
$ws.Range("E22").Value = "  -10.31%  "

$ws.Range("E23").Value = "  -1.36%  "

$ws.Range("D24").Value = "'15.95"
$ws.Range("E24").Value = "  -1.32%  "

$ws.Range("D25").Value = "'97.15"
$ws.Range("E25").Value = "  -0.85%  "

$ws.Range("E26").Value = "  -0.19%  "

$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  -0.04%  "

$ws.Range("E28").Value = "  -1.42%  "

$ws.Range("D29").Value = "'9.26"
$ws.Range("E29").Value = "  -1.56%  "

$ws.Range("D30").Value = "'33.47"
$ws.Range("E30").Value = "  +1.50%  "

$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'8.48"
$ws.Range("E31").Value = "  -2.75%  "

$ws.Range("B32").Value = "Stacks"
$ws.Range("C32").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D32").Value = "'3.11"
$ws.Range("E32").Value = "  -2.42%  "

$ws.Range("E33").Value = "  -1.30%  "

$ws.Range("D34").Value = "'7.05"
$ws.Range("E34").Value = "  -3.21%  "

$ws.Range("D35").Value = "'568.33"
$ws.Range("E35").Value = "  -4.90%  "

$ws.Range("D36").Value = "'3.67"
$ws.Range("E36").Value = "  +2.31%  "

$ws.Range("D37").Value = "'10.80"
$ws.Range("E37").Value = "  -0.69%  "

$ws.Range("D38").Value = "'57.45"
$ws.Range("E38").Value = "  +0.29%  "

$ws.Range("E39").Value = "  -1.73%  "

$ws.Range("E40").Value = "  -0.04%  "

$ws.Range("D41").Value = "'0.143"
$ws.Range("E41").Value = "  +5.70%  "

$ws.Range("E42").Value = "  +3.60%  "

$ws.Range("E43").Value = "  -1.88%  "

$ws.Range("D44").Value = "'3.337.25"
$ws.Range("E44").Value = "  -2.39%  "

$ws.Range("D45").Value = "'3.05"
$ws.Range("E45").Value = "  +4.96%  "

$ws.Range("D46").Value = "0.0₃0719"
$ws.Range("E46").Value = "  +1.70%  "

$ws.Range("D47").Value = "'33.14"
$ws.Range("E47").Value = "  -0.58%  "

$ws.Range("D48").Value = "'2.67"
$ws.Range("E48").Value = "  +2.43%  "

$ws.Range("E49").Value = "  -2.39%  "

$ws.Range("D50").Value = "'134.36"
$ws.Range("E50").Value = "  +0.91%  "

$ws.Range("D51").Value = "'5.74"
$ws.Range("E51").Value = "  -0.66%  "
